{"js": "// Remove the hyperlink that wraps the \"tidyverse\" mention right after\n// \"Make it easier on yourself by using \" \u2014 turn it back into plain\n// (still Courier-New / blue / underlined) text, same as Word's\n// \"Remove Hyperlink\" command. Other occurrences of \"tidyverse\" that are\n// not hyperlinks are left untouched.\nconst body = context.document.body;\nconst results = body.search(\"tidyverse\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const r = results.items[i];\n  r.load(\"hyperlink\");\n  await context.sync();\n\n  if (r.hyperlink) {\n    r.hyperlink = \"\";\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the hyperlink that wraps the \"tidyverse\" text right after\n# \"Make it easier on yourself by using \" and remove the hyperlink\n# (turn it back into plain, still-formatted text) while leaving the\n# other \"tidyverse\" mentions (which are not hyperlinks) untouched.\n$count = $d.Hyperlinks.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $h = $d.Hyperlinks.Item($i)\n    if ($h.TextToDisplay -eq \"tidyverse\") {\n        $h.Delete()\n    }\n}\n"}
